# SmartCutCalculation.xlsx reformat:
#  - header date / project name text updated
#  - the "Plank Length / Snippet [Floor/Apartment] / Waste" summary table gets
#    a "Columns" section title + a "6300 X 0" sub-header above it, and the
#    table headers are renamed to "Planks" / "Snippets [Floor/Apartment]" /
#    "Waste"; this pushes the whole data table down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three blank rows right before the old header row (row 4). This
#    shifts the entire existing table (old rows 4-58) down to rows 7-61,
#    which is exactly where it needs to land - no other row needs touching.
$ws.Range("4:6").Insert()

# 2) Update the title block (rows 1-2, untouched by the insert).
$ws.Range("C1").Value2 = "19.12.2019 08:27"
$ws.Range("A2").Value2 = "project name"

# 3) New "Columns" section title (row 4) - copy the bold header format from
#    the shifted-down table header row (row 7, style matches the original
#    A4/B4/C4 header).
$ws.Range("A7").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A4").Value2 = "Columns"

# 4) New "6300 X 0" sub-header (row 5) - copy the sub-header format from the
#    shifted-down "6900 X 27" cell (row 8, style matches the original A5).
$ws.Range("A8").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A5").Value2 = "6300 X 0"

# Row 6 stays blank (already is, from the insert).

# 5) Rename the table header cells (now on row 7).
$ws.Range("A7").Value2 = "Planks"
$ws.Range("B7").Value2 = "Snippets [Floor/Apartment]"
$ws.Range("C7").Value2 = "Waste"
